$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B/C/E/F are constant per block; column A (doctyp_code) and D (is_active)
# vary per row. Row 5 previously held "PSP" which is now replaced by the new
# "DOC###" document-type codes, and many additional rows are appended for the
# new Mac-Address / Document Types master data.

# Row 5: PSP -> DOC001 (category stays POI)
$ws.Cells.Item(5, 1).Value = "DOC001"

# Rows 8-19: doccat_code = POI (col B = 7 in old sheet, same text "POI")
$poiCodes = @("DOC001","DOC002","DOC003","DOC004","DOC005","DOC006","DOC007","DOC008","DOC009","DOC010","DOC011","DOC012")
for ($i = 0; $i -lt $poiCodes.Length; $i++) {
    $r = 8 + $i
    $ws.Cells.Item($r, 1).Value = $poiCodes[$i]
    $ws.Cells.Item($r, 2).Value = "POI"
    $ws.Cells.Item($r, 3).Value = "ara"
    $ws.Cells.Item($r, 4).Value = $true
    $ws.Cells.Item($r, 5).Value = "superadmin"
    $ws.Cells.Item($r, 6).Value = "now()"
}

# Rows 20-30: doccat_code = POA
$poaCodes = @("DOC001","DOC013","DOC014","DOC015","DOC004","DOC005","DOC006","DOC016","DOC017","DOC018","DOC008")
for ($i = 0; $i -lt $poaCodes.Length; $i++) {
    $r = 20 + $i
    $ws.Cells.Item($r, 1).Value = $poaCodes[$i]
    $ws.Cells.Item($r, 2).Value = "POA"
    $ws.Cells.Item($r, 3).Value = "ara"
    $ws.Cells.Item($r, 4).Value = $true
    $ws.Cells.Item($r, 5).Value = "superadmin"
    $ws.Cells.Item($r, 6).Value = "now()"
}

# Rows 31-36: doccat_code = CRN
$crnCodes = @("DOC024","DOC025","DOC026","DOC001","DOC027","DOC028")
for ($i = 0; $i -lt $crnCodes.Length; $i++) {
    $r = 31 + $i
    $ws.Cells.Item($r, 1).Value = $crnCodes[$i]
    $ws.Cells.Item($r, 2).Value = "CRN"
    $ws.Cells.Item($r, 3).Value = "ara"
    $ws.Cells.Item($r, 4).Value = $true
    $ws.Cells.Item($r, 5).Value = "superadmin"
    $ws.Cells.Item($r, 6).Value = "now()"
}

# Refresh selection to mirror the post-edit workbook state (cursor moved to
# column G, one column past the now-wider used range).
$ws.Range("G1:XFD1048576").Select()
